$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Text
    )
    $range = $ws.Range($Address)
    $range.NumberFormat = "@"
    $range.Value = $Text
    $range.Style = "Normal"
}

Set-TextValue 'D2' '35.624.49'
Set-TextValue 'E2' '  +0.60%  '
Set-TextValue 'D3' '1.900.12'
Set-TextValue 'E3' '  +0.42%  '
Set-TextValue 'E4' '  -0.07%  '
Set-TextValue 'B5' 'BNB'
Set-TextValue 'C5' 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
Set-TextValue 'D5' '247.01'
Set-TextValue 'E5' '  -0.06%  '
Set-TextValue 'B6' 'XRP'
Set-TextValue 'C6' 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
Set-TextValue 'D6' '0.691'
Set-TextValue 'E6' '  +0.14%  '
Set-TextValue 'E7' '  -0.02%  '
Set-TextValue 'D8' '43.33'
Set-TextValue 'E8' '  -1.93%  '
Set-TextValue 'E9' '  +1.71%  '
Set-TextValue 'D10' '56.35'
Set-TextValue 'E10' '  +8.37%  '
Set-TextValue 'E11' '  +2.80%  '
Set-TextValue 'E12' '  +1.51%  '
Set-TextValue 'D13' '14.39'
Set-TextValue 'E13' '  +9.54%  '
Set-TextValue 'D14' '0.799'
Set-TextValue 'E14' '  +10.21%  '
Set-TextValue 'D15' '2.175.72'
Set-TextValue 'E15' '  +0.31%  '
Set-TextValue 'E16' '  +1.98%  '
Set-TextValue 'D17' '1.900.34'
Set-TextValue 'E17' '  -0.13%  '
Set-TextValue 'D18' '35.605.13'
Set-TextValue 'E18' '  +0.60%  '
Set-TextValue 'D19' '73.79'
Set-TextValue 'E19' '  +0.68%  '
Set-TextValue 'E20' '  +1.50%  '
Set-TextValue 'D21' '246.29'
Set-TextValue 'E21' '  +0.04%  '
Set-TextValue 'E22' '  +1.69%  '
Set-TextValue 'D23' '5.22'
Set-TextValue 'E23' '  +5.16%  '
Set-TextValue 'E24' '  +4.83%  '
Set-TextValue 'E25' '  -0.12%  '
Set-TextValue 'D26' '2.17'
Set-TextValue 'E26' '  -1.25%  '
Set-TextValue 'D27' '167.01'
Set-TextValue 'E27' '  +0.88%  '
Set-TextValue 'D28' '8.65'
Set-TextValue 'E28' '  +1.74%  '
Set-TextValue 'D29' '18.38'
Set-TextValue 'E29' '  +0.21%  '
Set-TextValue 'E30' '  +0.54%  '
Set-TextValue 'D31' '4.38'
Set-TextValue 'E31' '  +2.77%  '
Set-TextValue 'D32' '0.0604'
Set-TextValue 'E32' '  +3.80%  '
Set-TextValue 'D33' '4.26'
Set-TextValue 'E33' '  +0.32%  '
Set-TextValue 'E34' '  +16.12%  '
Set-TextValue 'D36' '1.49'
Set-TextValue 'E36' '  -16.89%  '
Set-TextValue 'D37' '0.856'
Set-TextValue 'E37' '  +0.50%  '
Set-TextValue 'B38' 'LidoDAOToken'
Set-TextValue 'C38' 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue 'D38' '1.96'
Set-TextValue 'E38' '  -1.99%  '
Set-TextValue 'B39' 'Kaspa'
Set-TextValue 'C39' 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D39' '0.0739'
Set-TextValue 'E39' '  +10.23%  '
Set-TextValue 'D40' '0.0226'
Set-TextValue 'E40' '  +5.88%  '
Set-TextValue 'D41' '99.62'
Set-TextValue 'E41' '  +2.05%  '
Set-TextValue 'D42' '17.12'
Set-TextValue 'E42' '  -0.66%  '
Set-TextValue 'E43' '  -1.03%  '
Set-TextValue 'D44' '13.60'
Set-TextValue 'E44' '  +12.58%  '
Set-TextValue 'D45' '1.324.54'
Set-TextValue 'E45' '  +2.72%  '
Set-TextValue 'E46' '  +0.47%  '
Set-TextValue 'D47' '0.0813'
Set-TextValue 'E47' '  +1.07%  '
Set-TextValue 'E48' '  +0.79%  '
Set-TextValue 'E49' '  -0.34%  '
Set-TextValue 'D50' '6.39'
Set-TextValue 'E50' '  -0.38%  '
Set-TextValue 'D51' '42.62'
Set-TextValue 'E51' '  -1.03%  '
